# Update "想去人数" (F column) values on the 展览 and 全部类型 sheets.
# The same set of rows/values need to change identically on both sheets.

$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 8878
    3  = 8337
    4  = 151
    5  = 164
    7  = 252
    8  = 756
    9  = 214
    10 = 5451
    13 = 90
    14 = 19
    17 = 162
    18 = 203
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
